# Append a new user record (Cliff Fritsch) to the userData sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data starts on row 2; find the first empty row right after the existing block.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "Cliff"
$ws.Cells.Item($newRow, 2).Value = "Fritsch"
$ws.Cells.Item($newRow, 3).Value = "alberto.gulgowski@example.com"
$ws.Cells.Item($newRow, 4).Value = "DO5LaP%mZ"

# The date-of-birth column holds plain text like "MM/DD/YYYY" in this sheet,
# not a real date value, so force text formatting before assigning it -
# otherwise Excel would silently convert the string into a date serial
# number. Clear the formatting afterwards so the cell ends up with the
# same (default) style as every other cell in the sheet.
$ws.Cells.Item($newRow, 5).NumberFormat = "@"
$ws.Cells.Item($newRow, 5).Value = "04/16/1979"
$ws.Cells.Item($newRow, 5).ClearFormats()
